# Insert a new row at row 932 (pushes existing rows 932..1032 down to 933..1033)
# and populate the new row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 932 (limited to the used columns A:T so
# that Excel does not widen the row's column span to the full 16384 columns);
# this shifts the existing rows 932..1032 down to 933..1033.
$ws.Range("A932:T932").Insert(-4121)  # xlShiftDown

# Copy the formatting of the row above (row 931, which itself mirrors the prior
# data rows) down onto the freshly inserted row 932 so styles (e.g. date format
# in column D) match the rest of the table.
$ws.Range("A931:T931").Copy()
$ws.Range("A932:T932").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 932 with the new data record.
$ws.Cells.Item(932, 1).Value = 7
$ws.Cells.Item(932, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(932, 3).Value = "Ñuble"
$ws.Cells.Item(932, 4).Value = 45132
$ws.Cells.Item(932, 5).Value = 16
$ws.Cells.Item(932, 6).Value = "Fruta"
$ws.Cells.Item(932, 7).Value = 100106
$ws.Cells.Item(932, 8).Value = "Oleaginosos"
$ws.Cells.Item(932, 9).Value = 100106002
$ws.Cells.Item(932, 10).Value = "Palta"
$ws.Cells.Item(932, 11).Value = "Hass"
$ws.Cells.Item(932, 12).Value = "Primera"
$ws.Cells.Item(932, 13).Value = 120
$ws.Cells.Item(932, 14).Value = 24000
$ws.Cells.Item(932, 15).Value = 25000
$ws.Cells.Item(932, 16).Value = 24500
$ws.Cells.Item(932, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(932, 18).Value = "Perú"
$ws.Cells.Item(932, 19).Value = 2450
$ws.Cells.Item(932, 20).Value = 10
